# Update data: 2025-11-05 13:28
$wb = $excel.ActiveWorkbook

# 1. Update "Last Updated" timestamp on the Metadata sheet (A2)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 01:28 PM"

# 2. Update the "1 Year" column (F) on the Industry Analysis sheet, rows 2-76
$ind = $wb.Worksheets.Item("Industry Analysis")

$newValues = @(
    21.0016,-16.2396,27.1317,-50.6494,53.2813,-8.106199999999999,-9.552099999999999,
    36.3756,-6.1314,31.9081,-18.4955,14.0155,-36.0718,-0.1622,0.1459,-22.0012,1.0561,
    -27.708,47.7309,12.0959,95.1491,-50.2657,-13.3427,-9.9316,5.8244,-32.7692,-24.8224,
    -18.4191,25.8569,58.4712,-3.3862,-6.3282,27.7203,4.4873,-4.9458,3.6074,-23.3973,
    8.7355,-5.8541,-8.3934,20.3818,14.3164,-12.6846,28.4075,-1.1135,-37.1997,-29.8569,
    -27.5511,-49.7478,-51.8002,-38.5254,-12.4886,-5.0725,-17.7445,-26.636,-29.3361,
    -11.9574,-24.5687,-12.3,-10.9446,-17.1229,-9.5038,54.2749,-43.4736,13.2687,12.7149,
    24.8057,-17.0328,-6.8927,13.6034,3.9995,-16.226,-16.2448,28.6924,48.9752
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ind.Cells.Item($row, 6).Value = $newValues[$i]
}
